{"js": "// Replace the 100 math-problem strings in the 20x5 \"within 100\" table,\n// in document order, with their updated versions (per the commit diff).\n// Each table cell holds exactly one run of text such as \"93-23=\" and\n// the diff changes every one of those strings (keeping the surrounding\n// paragraph/run formatting exactly as-is).\n\nconst newValues = [\n  [\"94-90=\", \"23+75=\", \"99-37=\", \"47+37=\", \"33+46=\"],\n  [\"6+6=\", \"81-63=\", \"1+93=\", \"54+6=\", \"90-16=\"],\n  [\"42+55=\", \"73+3=\", \"44-29=\", \"77+13=\", \"50-30=\"],\n  [\"86-20=\", \"30-18=\", \"30+40=\", \"24+68=\", \"2+6=\"],\n  [\"36-20=\", \"46+22=\", \"4+73=\", \"90-15=\", \"75-42=\"],\n  [\"30+52=\", \"39-39=\", \"95-11=\", \"65+1=\", \"46+46=\"],\n  [\"21+48=\", \"87-73=\", \"46-25=\", \"58-29=\", \"98-64=\"],\n  [\"45-16=\", \"63-31=\", \"69-17=\", \"63+7=\", \"86-66=\"],\n  [\"2+80=\", \"92-14=\", \"98-49=\", \"86-74=\", \"44+1=\"],\n  [\"11+56=\", \"37-7=\", \"65+26=\", \"25+40=\", \"31+1=\"],\n  [\"66+27=\", \"24+7=\", \"24+66=\", \"49-4=\", \"77-76=\"],\n  [\"24+2=\", \"62-19=\", \"10+6=\", \"19+79=\", \"35+31=\"],\n  [\"13+47=\", \"92-32=\", \"15+44=\", \"49-32=\", \"71+20=\"],\n  [\"58-12=\", \"38+49=\", \"24+10=\", \"86-44=\", \"71-61=\"],\n  [\"79-70=\", \"21-1=\", \"70-12=\", \"17+43=\", \"16+80=\"],\n  [\"66+8=\", \"53+22=\", \"99-5=\", \"3-0=\", \"5+33=\"],\n  [\"91-47=\", \"95-23=\", \"50-18=\", \"85-20=\", \"49-36=\"],\n  [\"60+33=\", \"4+80=\", \"21+8=\", \"15+76=\", \"45+23=\"],\n  [\"25-6=\", \"7+80=\", \"43+40=\", \"84-73=\", \"73-58=\"],\n  [\"55-28=\", \"97-2=\", \"73+26=\", \"67-53=\", \"57+24=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 math-problem strings in the 20x5 \"within 100\" table,\n# in document (row-major) order, with their updated versions per the commit diff.\n# Each table cell holds exactly one run of text such as \"93-23=\" and every one\n# of those strings changes; surrounding paragraph/run formatting is untouched\n# because we only rewrite Cell(r,c).Range.Text (which excludes the end-of-cell mark).\n\n$newValues = @(\n    \"94-90=\",\n    \"23+75=\",\n    \"99-37=\",\n    \"47+37=\",\n    \"33+46=\",\n    \"6+6=\",\n    \"81-63=\",\n    \"1+93=\",\n    \"54+6=\",\n    \"90-16=\",\n    \"42+55=\",\n    \"73+3=\",\n    \"44-29=\",\n    \"77+13=\",\n    \"50-30=\",\n    \"86-20=\",\n    \"30-18=\",\n    \"30+40=\",\n    \"24+68=\",\n    \"2+6=\",\n    \"36-20=\",\n    \"46+22=\",\n    \"4+73=\",\n    \"90-15=\",\n    \"75-42=\",\n    \"30+52=\",\n    \"39-39=\",\n    \"95-11=\",\n    \"65+1=\",\n    \"46+46=\",\n    \"21+48=\",\n    \"87-73=\",\n    \"46-25=\",\n    \"58-29=\",\n    \"98-64=\",\n    \"45-16=\",\n    \"63-31=\",\n    \"69-17=\",\n    \"63+7=\",\n    \"86-66=\",\n    \"2+80=\",\n    \"92-14=\",\n    \"98-49=\",\n    \"86-74=\",\n    \"44+1=\",\n    \"11+56=\",\n    \"37-7=\",\n    \"65+26=\",\n    \"25+40=\",\n    \"31+1=\",\n    \"66+27=\",\n    \"24+7=\",\n    \"24+66=\",\n    \"49-4=\",\n    \"77-76=\",\n    \"24+2=\",\n    \"62-19=\",\n    \"10+6=\",\n    \"19+79=\",\n    \"35+31=\",\n    \"13+47=\",\n    \"92-32=\",\n    \"15+44=\",\n    \"49-32=\",\n    \"71+20=\",\n    \"58-12=\",\n    \"38+49=\",\n    \"24+10=\",\n    \"86-44=\",\n    \"71-61=\",\n    \"79-70=\",\n    \"21-1=\",\n    \"70-12=\",\n    \"17+43=\",\n    \"16+80=\",\n    \"66+8=\",\n    \"53+22=\",\n    \"99-5=\",\n    \"3-0=\",\n    \"5+33=\",\n    \"91-47=\",\n    \"95-23=\",\n    \"50-18=\",\n    \"85-20=\",\n    \"49-36=\",\n    \"60+33=\",\n    \"4+80=\",\n    \"21+8=\",\n    \"15+76=\",\n    \"45+23=\",\n    \"25-6=\",\n    \"7+80=\",\n    \"43+40=\",\n    \"84-73=\",\n    \"73-58=\",\n    \"55-28=\",\n    \"97-2=\",\n    \"73+26=\",\n    \"67-53=\",\n    \"57+24=\",\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n\n"}
